$d = $word.ActiveDocument

# --- 1. Split off a blank paragraph right after "Project Repository" ---
$d.Content.Find.Execute("Project Repository", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Project Repository^p", 2) | Out-Null

$blank = $d.Paragraphs.Item(13)
$blank.Range.ListFormat.RemoveNumbers()
$blank.Style = "Normal"

# --- 2. Make room for the new hyperlink paragraph right after the bookmark paragraph,
#        before we touch its text/formatting (so the new paragraph stays unformatted) ---
$bmPara = $d.Paragraphs.Item(14)
$bmPara.Range.InsertParagraphAfter()

$linkPara = $d.Paragraphs.Item(15)
$d.Hyperlinks.Add($linkPara.Range, "https://github.com/pholohan/GreenGuard", "", "", `
                   "https://github.com/pholohan/GreenGuard") | Out-Null

# --- 3. Turn the bookmarked paragraph's text into the bold "GIT Repository" label,
#        keeping the _GoBack bookmark that already lives in that paragraph ---
$bmPara = $d.Paragraphs.Item(14)
$bmPara.Range.Text = "GIT Repository"
$bmPara.Range.Bold = 1

# --- 4. Define the (now referenced) Hyperlink character style ---
$hlStyle = $d.Styles.Add("Hyperlink", 2)
$hlStyle.BaseStyle = "DefaultParagraphFont"
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = 1
$hlStyle.Font.Color = 16711680
$hlStyle.Font.Underline = 1

Write-Host "done"
